$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation (row) was recorded for Espinaca at
# "Vega Monumental Concepción". It belongs right after the existing
# most-recent entry (old row 83) and before the rest of the historical
# series, so insert a fresh row at 84 and push everything else down by one.
$ws.Rows.Item(84).Insert()

$ws.Range("A84").Value2 = 11
$ws.Range("B84").Value2 = "Vega Monumental Concepción"
$ws.Range("C84").Value2 = "Bíobío"
$ws.Range("D84").Value2 = 45141
$ws.Range("E84").Value2 = 8
$ws.Range("F84").Value2 = 100112012
$ws.Range("G84").Value2 = "Espinaca"
$ws.Range("H84").Value2 = "Sin especificar"
$ws.Range("I84").Value2 = "Primera"
$ws.Range("J84").Value2 = 70
$ws.Range("K84").Value2 = 6000
$ws.Range("L84").Value2 = 6500
$ws.Range("M84").Value2 = 6143
$ws.Range("N84").Value2 = "$/cuna 10 kilos"
$ws.Range("O84").Value2 = "Región Metropolitana"
$ws.Range("P84").Value2 = 614
$ws.Range("Q84").Value2 = 10
$ws.Range("R84").Value2 = "Hortaliza"
